# Testing github action (xlsx)
# Replace the Korean "Invalid log path" translation (LogPathInvalid / ko-KR, cell C9)
# with a revised string that is entered as rich text: the Korean words use the
# "돋움" font while the ASCII spaces/period stay in Calibri - matching how the
# text was retyped character-by-character in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Msgs")

$cell = $ws.Range("C9")
$cell.Value = "로그 경로가 유효합니다."

# Rich-text runs (1-based character positions into the new string):
#   1-2   "로그"      -> 돋움 10pt
#   3     " "         -> Calibri 10pt
#   4-6   "경로가"    -> 돋움 10pt
#   7     " "         -> Calibri 10pt
#   8-12  "유효합니다" -> 돋움 10pt
#   13    "."         -> Calibri 10pt

$run1 = $cell.Characters(1, 2)
$run1.Font.Name = "돋움"
$run1.Font.Size = 10

$run2 = $cell.Characters(3, 1)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 10

$run3 = $cell.Characters(4, 3)
$run3.Font.Name = "돋움"
$run3.Font.Size = 10

$run4 = $cell.Characters(7, 1)
$run4.Font.Name = "Calibri"
$run4.Font.Size = 10

$run5 = $cell.Characters(8, 5)
$run5.Font.Name = "돋움"
$run5.Font.Size = 10

$run6 = $cell.Characters(13, 1)
$run6.Font.Name = "Calibri"
$run6.Font.Size = 10

# Leave the cell selected on C9, matching the post-edit UI state.
$cell.Select()

Write-Output "Updated C9 rich text"
